$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("missing_stations")

# Sort the "free" block (A2:G41) ascending by column A (id), mirroring the
# manual Data > Sort the author performed while consolidating the table.
$rng = $ws.Range("A2:G41")
$key1 = $ws.Range("A2:A41")
$rng.Sort($key1, 1, $null, $null, 1, $null, 1, 2)

# Manually swap the two untagged (blank id) rows that follow - the author
# reordered "DIVVY Map Frame B/C Station" and "DIVVY CASSETTE REPAIR MOBILE
# STATION" by hand.
$b42 = $ws.Range("B42").Value2
$b43 = $ws.Range("B43").Value2
$ws.Range("B42").Value = $b43
$ws.Range("B43").Value = $b42

# Update the view state to match where the author ended up working.
$ws.Range("A18").Select()
$av = $excel.ActiveWindow
$av.FreezePanes = $false
$ws.Range("A2").Select()
$av.FreezePanes = $true

$ws.Range("D27").Select()
